$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Manual update of the H index value in A3 (commit: "H index manual update")
$ws.Range("A3").Value = 585

# Update the active cell selection to A4, matching the saved selection state
$null = $ws.Range("A4").Select()
